# Updates cryptos list: refresh Price (D) and Volume(1h) (E) values for
# each coin row, and re-sort three coin pairs whose rank order changed
# (Filecoin/ImmutableX, HuobiToken/VeChain, Aave/FraxShare) by swapping
# their Coin (B), Link (C), Price (D) and Volume (E) cell contents.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '36.176.42'
$ws.Range("E2").Value = '  -1.11%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.015.26'
$ws.Range("E3").Value = '  -1.41%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.96'
$ws.Range("E5").Value = '  +3.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.643'
$ws.Range("E6").Value = '  -3.51%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.63'
$ws.Range("E7").Value = '  +12.55%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '59.03'
$ws.Range("E9").Value = '  -6.70%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.372'
$ws.Range("E10").Value = '  +1.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0750'
$ws.Range("E11").Value = '  +0.38%  '
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.919'
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.00'
$ws.Range("E14").Value = '  +3.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.307.17'
$ws.Range("E15").Value = '  -1.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.43'
$ws.Range("E16").Value = '  +0.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '19.56'
$ws.Range("E17").Value = '  +13.14%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.036.80'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '36.103.00'
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.21'
$ws.Range("E20").Value = '  +0.26%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.0₃0861'
$ws.Range("E21").Value = '  +0.61%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.28'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '234.22'
$ws.Range("E23").Value = '  -1.59%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.73'
$ws.Range("E24").Value = '  +21.44%  '
$ws.Range("E25").Value = '  +0.21%  '
$ws.Range("E26").Value = '  -2.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.55'
$ws.Range("E27").Value = '  +3.88%  '
$ws.Range("E28").Value = '  +0.27%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.65'
$ws.Range("E29").Value = '  -1.64%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.121'
$ws.Range("E30").Value = '  -0.54%  '
$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.22'
$ws.Range("E31").Value = '  +2.98%  '
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.15'
$ws.Range("E32").Value = '  +3.15%  '
$ws.Range("E33").Value = '  +24.60%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0606'
$ws.Range("E34").Value = '  +1.05%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.52'
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.49'
$ws.Range("E36").Value = '  +13.87%  '
$ws.Range("E37").Value = '  -0.19%  '
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("E39").Value = '  +16.14%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.105'
$ws.Range("E40").Value = '  +16.67%  '
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("B42").Value = 'VeChain'
$ws.Range("C42").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0216'
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.89'
$ws.Range("E43").Value = '  +1.29%  '
$ws.Range("E44").Value = '  +1.97%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '16.80'
$ws.Range("E45").Value = '  +5.74%  '
$ws.Range("B46").Value = 'FraxShare'
$ws.Range("C46").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.85'
$ws.Range("E46").Value = '  +4.61%  '
$ws.Range("B47").Value = 'Aave'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '94.08'
$ws.Range("E47").Value = '  +0.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.419.78'
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.43'
$ws.Range("E49").Value = '  +8.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.90'
$ws.Range("E50").Value = '  -1.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '47.30'
$ws.Range("E51").Value = '  +3.51%  '

